# Update Maggie Burton's "Time Spent" value on the timecards sheet
# from "8h 30m" to "14h 30m".
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B7").Value = "14h 30m"
